$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new metric data row (row 17) with timestamp text and a numeric value
$ws.Range("A17").Value = "2025-04-28 20:53:44"
$ws.Range("B17").Value = 0
